$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restyle the surviving "2014" header cell (D5) so it carries the same
# border/format that the "1989" header cell (B5) had, before we collapse
# the 1989/2002/2014 columns down to a single remaining column.
$ws.Cells.Item(5,2).Copy()
$ws.Cells.Item(5,4).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Drop the "(according to the population census data)" sub-title row.
$ws.Rows(2).Delete()

# Drop the now-obsolete "1989" and "2002" columns, leaving only "2014".
$ws.Columns("B:C").Delete()

# Row heights: the remaining rows grow from the compact report layout to
# the taller 20.1pt rows used in the trimmed-down sheet.
$ws.Rows(1).RowHeight = 20.1
$ws.Rows(2).RowHeight = 20.1
$ws.Rows(3).RowHeight = 20.1
$ws.Rows(4).RowHeight = 20.1
$ws.Rows(5).RowHeight = 20.1
